$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 10503
$ws.Range("F3").Value = 436
$ws.Range("F6").Value = 287
$ws.Range("F9").Value = 783
$ws.Range("F12").Value = 1111
$ws.Range("F13").Value = 3258
$ws.Range("F14").Value = 2422
$ws.Range("F16").Value = 2180
$ws.Range("F17").Value = 2180
$ws.Range("F19").Value = 1942
$ws.Range("F21").Value = 1597
$ws.Range("F22").Value = 588
$ws.Range("F24").Value = 259
$ws.Range("F25").Value = 11
$ws.Range("F26").Value = 28
$ws.Range("F29").Value = 390
$ws.Range("F31").Value = 78
$ws.Range("F32").Value = 408
$ws.Range("F33").Value = 608
$ws.Range("F34").Value = 32
$ws.Range("F35").Value = 57
$ws.Range("F36").Value = 274
$ws.Range("F37").Value = 12
$ws.Range("F39").Value = 507
$ws.Range("F40").Value = 483
$ws.Range("F41").Value = 1732
$ws.Range("F42").Value = 146
$ws.Range("F43").Value = 451
$ws.Range("F44").Value = 58
$ws.Range("F45").Value = 475
$ws.Range("F46").Value = 1054
$ws.Range("F48").Value = 368

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 30
$ws.Range("F4").Value = 48
$ws.Range("F8").Value = 1

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 10503
$ws.Range("F3").Value = 436
$ws.Range("F5").Value = 30
$ws.Range("F8").Value = 287
$ws.Range("F11").Value = 783
$ws.Range("F12").Value = 1111
$ws.Range("F13").Value = 3258
$ws.Range("F14").Value = 2422
$ws.Range("F15").Value = 2180
$ws.Range("F16").Value = 2180
$ws.Range("F17").Value = 1597
$ws.Range("F18").Value = 588
$ws.Range("F20").Value = 259
$ws.Range("F21").Value = 11
$ws.Range("F22").Value = 28
$ws.Range("F25").Value = 390
$ws.Range("F27").Value = 78
$ws.Range("F28").Value = 408
$ws.Range("F29").Value = 608
$ws.Range("F30").Value = 32
$ws.Range("F31").Value = 48
$ws.Range("F34").Value = 57
$ws.Range("F35").Value = 274
$ws.Range("F37").Value = 507
$ws.Range("F39").Value = 483
$ws.Range("F40").Value = 1732
$ws.Range("F41").Value = 146
$ws.Range("F42").Value = 1
$ws.Range("F45").Value = 451
$ws.Range("F46").Value = 58
$ws.Range("F47").Value = 475
$ws.Range("F48").Value = 1054
$ws.Range("F49").Value = 368
